# Insert a new column before column D (shifts existing D:K data to E:L)
# and backfill column D plus the updated values in the quarterly table with
# the refreshed financial figures (VMW quarterly financials update).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a blank column at D; existing D:K shift right to E:L.
$ws.Columns("D").Insert()

# 2. Copy the number/date formatting from the (now shifted) column E onto the
#    new column D so the new quarter's cells match the existing look
#    (date format for the period-ending rows, #,##0 for the data rows).
#    Done per-table so the blank separator/header rows (36/37, 78/79) that
#    never had a column-D cell to begin with stay untouched.
$ws.Range("E7:E35").Copy()
$ws.Range("D7:D35").PasteSpecial(-4122)
$ws.Range("E38:E77").Copy()
$ws.Range("D38:D77").PasteSpecial(-4122)
$ws.Range("E80:E102").Copy()
$ws.Range("D80:D102").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false



# 3. Write the full refreshed data block (columns D:L) for each of the three
#    tables (Income Statement rows 7-35, Balance Sheet rows 38-77, Cash Flow
#    Statement rows 80-102).
$b1 = New-Object 'object[,]' 29,9
$b1[0,0] = 43406
$b1[0,1] = 43315
$b1[0,2] = 43224
$b1[0,3] = 43133
$b1[0,4] = 43042
$b1[0,5] = 42951
$b1[0,6] = 42860
$b1[0,7] = 42735
$b1[0,8] = 42643
$b1[1,0] = 2200000
$b1[1,1] = 2174000
$b1[1,2] = 2008000
$b1[1,3] = 2226000
$b1[1,4] = 1938000
$b1[1,5] = 1932000
$b1[1,6] = 1765000
$b1[1,7] = 2033000
$b1[1,8] = 1778000
$b1[2,0] = 315000
$b1[2,1] = 305000
$b1[2,2] = 296000
$b1[2,3] = 305000
$b1[2,4] = 278000
$b1[2,5] = 270000
$b1[2,6] = 289000
$b1[2,7] = 274000
$b1[2,8] = 266000
$b1[3,0] = 1885000
$b1[3,1] = 1869000
$b1[3,2] = 1712000
$b1[3,3] = 1921000
$b1[3,4] = 1660000
$b1[3,5] = 1662000
$b1[3,6] = 1476000
$b1[3,7] = 1759000
$b1[3,8] = 1512000
$b1[4,0] = $null
$b1[4,1] = $null
$b1[4,2] = $null
$b1[4,3] = $null
$b1[4,4] = $null
$b1[4,5] = $null
$b1[4,6] = $null
$b1[4,7] = $null
$b1[4,8] = $null
$b1[5,0] = 499000
$b1[5,1] = 481000
$b1[5,2] = 453000
$b1[5,3] = 456000
$b1[5,4] = 449000
$b1[5,5] = 428000
$b1[5,6] = 421000
$b1[5,7] = 394000
$b1[5,8] = 389000
$b1[6,0] = 0
$b1[6,1] = 0
$b1[6,2] = 0
$b1[6,3] = 0
$b1[6,4] = 0
$b1[6,5] = 0
$b1[6,6] = 0
$b1[6,7] = 0
$b1[6,8] = 0
$b1[7,0] = 6000
$b1[7,1] = 1000
$b1[7,2] = 2000
$b1[7,3] = 10000
$b1[7,4] = 2000
$b1[7,5] = 36000
$b1[7,6] = 64000
$b1[7,7] = 0
$b1[7,8] = 0
$b1[8,0] = "NA"
$b1[8,1] = "NA"
$b1[8,2] = "NA"
$b1[8,3] = 7000
$b1[8,4] = "NA"
$b1[8,5] = "NA"
$b1[8,6] = "NA"
$b1[8,7] = 0
$b1[8,8] = 0
$b1[9,0] = $null
$b1[9,1] = $null
$b1[9,2] = $null
$b1[9,3] = $null
$b1[9,4] = $null
$b1[9,5] = $null
$b1[9,6] = $null
$b1[9,7] = $null
$b1[9,8] = $null
$b1[10,0] = 1705000
$b1[10,1] = 1665000
$b1[10,2] = 1626000
$b1[10,3] = 1621000
$b1[10,4] = 1528000
$b1[10,5] = 1507000
$b1[10,6] = 1504000
$b1[10,7] = 1490000
$b1[10,8] = 1397000
$b1[11,0] = 495000
$b1[11,1] = 509000
$b1[11,2] = 382000
$b1[11,3] = 605000
$b1[11,4] = 410000
$b1[11,5] = 425000
$b1[11,6] = 261000
$b1[11,7] = 543000
$b1[11,8] = 381000
$b1[12,0] = $null
$b1[12,1] = $null
$b1[12,2] = $null
$b1[12,3] = $null
$b1[12,4] = $null
$b1[12,5] = $null
$b1[12,6] = $null
$b1[12,7] = $null
$b1[12,8] = $null
$b1[13,0] = -117000
$b1[13,1] = 297000
$b1[13,2] = 827000
$b1[13,3] = 53000
$b1[13,4] = 31000
$b1[13,5] = 76000
$b1[13,6] = 27000
$b1[13,7] = 12000
$b1[13,8] = 13000
$b1[14,0] = 536000
$b1[14,1] = 957000
$b1[14,2] = 1365000
$b1[14,3] = 804000
$b1[14,4] = 590000
$b1[14,5] = 635000
$b1[14,6] = 424000
$b1[14,7] = 639000
$b1[14,8] = 481000
$b1[15,0] = 33000
$b1[15,1] = 34000
$b1[15,2] = 34000
$b1[15,3] = 33000
$b1[15,4] = 28000
$b1[15,5] = 7000
$b1[15,6] = 7000
$b1[15,7] = 6000
$b1[15,8] = 7000
$b1[16,0] = 345000
$b1[16,1] = 772000
$b1[16,2] = 1175000
$b1[16,3] = 625000
$b1[16,4] = 413000
$b1[16,5] = 494000
$b1[16,6] = 281000
$b1[16,7] = 549000
$b1[16,8] = 387000
$b1[17,0] = 11000
$b1[17,1] = 128000
$b1[17,2] = 233000
$b1[17,3] = 42000
$b1[17,4] = 18000
$b1[17,5] = 88000
$b1[17,6] = 36000
$b1[17,7] = 108000
$b1[17,8] = 68000
$b1[18,0] = 0
$b1[18,1] = 0
$b1[18,2] = 0
$b1[18,3] = 0
$b1[18,4] = 0
$b1[18,5] = 0
$b1[18,6] = 0
$b1[18,7] = 0
$b1[18,8] = 0
$b1[19,0] = 334000
$b1[19,1] = 644000
$b1[19,2] = 942000
$b1[19,3] = 583000
$b1[19,4] = 395000
$b1[19,5] = 406000
$b1[19,6] = 245000
$b1[19,7] = 441000
$b1[19,8] = 319000
$b1[20,0] = 334000
$b1[20,1] = 644000
$b1[20,2] = 942000
$b1[20,3] = 583000
$b1[20,4] = 395000
$b1[20,5] = 406000
$b1[20,6] = 245000
$b1[20,7] = 441000
$b1[20,8] = 319000
$b1[21,0] = 0
$b1[21,1] = 0
$b1[21,2] = 0
$b1[21,3] = 0
$b1[21,4] = 0
$b1[21,5] = 0
$b1[21,6] = 0
$b1[21,7] = 0
$b1[21,8] = 0
$b1[22,0] = "NA"
$b1[22,1] = "NA"
$b1[22,2] = "NA"
$b1[22,3] = -970000
$b1[22,4] = "NA"
$b1[22,5] = "NA"
$b1[22,6] = "NA"
$b1[22,7] = "NA"
$b1[22,8] = "NA"
$b1[23,0] = 0
$b1[23,1] = 0
$b1[23,2] = 0
$b1[23,3] = 0
$b1[23,4] = 0
$b1[23,5] = 0
$b1[23,6] = 0
$b1[23,7] = 0
$b1[23,8] = 0
$b1[24,0] = 0
$b1[24,1] = 0
$b1[24,2] = 0
$b1[24,3] = 0
$b1[24,4] = 0
$b1[24,5] = 0
$b1[24,6] = 0
$b1[24,7] = 0
$b1[24,8] = 0
$b1[25,0] = 117000
$b1[25,1] = -297000
$b1[25,2] = -827000
$b1[25,3] = -53000
$b1[25,4] = -31000
$b1[25,5] = -76000
$b1[25,6] = -27000
$b1[25,7] = -12000
$b1[25,8] = -13000
$b1[26,0] = 334000
$b1[26,1] = 644000
$b1[26,2] = 942000
$b1[26,3] = -387000
$b1[26,4] = 395000
$b1[26,5] = 406000
$b1[26,6] = 245000
$b1[26,7] = 441000
$b1[26,8] = 319000
$b1[27,0] = 0
$b1[27,1] = 0
$b1[27,2] = 0
$b1[27,3] = 0
$b1[27,4] = 0
$b1[27,5] = 0
$b1[27,6] = 0
$b1[27,7] = 0
$b1[27,8] = 0
$b1[28,0] = 334000
$b1[28,1] = 644000
$b1[28,2] = 942000
$b1[28,3] = -387000
$b1[28,4] = 395000
$b1[28,5] = 406000
$b1[28,6] = 245000
$b1[28,7] = 441000
$b1[28,8] = 319000
$ws.Range("D7:L35").Value = $b1

$b2 = New-Object 'object[,]' 40,9
$b2[0,0] = 43406
$b2[0,1] = 43315
$b2[0,2] = 43224
$b2[0,3] = 43133
$b2[0,4] = 43042
$b2[0,5] = 42951
$b2[0,6] = 42860
$b2[0,7] = 42735
$b2[0,8] = 42643
$b2[1,0] = $null
$b2[1,1] = $null
$b2[1,2] = $null
$b2[1,3] = $null
$b2[1,4] = $null
$b2[1,5] = $null
$b2[1,6] = $null
$b2[1,7] = $null
$b2[1,8] = $null
$b2[2,0] = $null
$b2[2,1] = $null
$b2[2,2] = $null
$b2[2,3] = $null
$b2[2,4] = $null
$b2[2,5] = $null
$b2[2,6] = $null
$b2[2,7] = $null
$b2[2,8] = $null
$b2[3,0] = 9189000
$b2[3,1] = 8121000
$b2[3,2] = 7101000
$b2[3,3] = 5971000
$b2[3,4] = 6012000
$b2[3,5] = 3552000
$b2[3,6] = 3864000
$b2[3,7] = 2790000
$b2[3,8] = 2654000
$b2[4,0] = 4338000
$b2[4,1] = 5179000
$b2[4,2] = 5529000
$b2[4,3] = 5682000
$b2[4,4] = 5600000
$b2[4,5] = 5350000
$b2[4,6] = 4748000
$b2[4,7] = 5195000
$b2[4,8] = 5600000
$b2[5,0] = 1664000
$b2[5,1] = 1729000
$b2[5,2] = 1202000
$b2[5,3] = 1953000
$b2[5,4] = 1154000
$b2[5,5] = 1396000
$b2[5,6] = 994000
$b2[5,7] = 1988000
$b2[5,8] = 1142000
$b2[6,0] = 0
$b2[6,1] = 0
$b2[6,2] = 0
$b2[6,3] = 0
$b2[6,4] = 0
$b2[6,5] = 0
$b2[6,6] = 0
$b2[6,7] = 0
$b2[6,8] = 0
$b2[7,0] = 227000
$b2[7,1] = 219000
$b2[7,2] = 269000
$b2[7,3] = 467000
$b2[7,4] = 160000
$b2[7,5] = 173000
$b2[7,6] = 172000
$b2[7,7] = 362000
$b2[7,8] = 159000
$b2[8,0] = 15418000
$b2[8,1] = 15248000
$b2[8,2] = 14101000
$b2[8,3] = 13836000
$b2[8,4] = 12926000
$b2[8,5] = 10471000
$b2[8,6] = 9778000
$b2[8,7] = 10335000
$b2[8,8] = 9555000
$b2[9,0] = 0
$b2[9,1] = 0
$b2[9,2] = 0
$b2[9,3] = 0
$b2[9,4] = 0
$b2[9,5] = 0
$b2[9,6] = 0
$b2[9,7] = 0
$b2[9,8] = 0
$b2[10,0] = 1128000
$b2[10,1] = 1105000
$b2[10,2] = 1098000
$b2[10,3] = 1074000
$b2[10,4] = 1031000
$b2[10,5] = 1005000
$b2[10,6] = 993000
$b2[10,7] = 1049000
$b2[10,8] = 1050000
$b2[11,0] = 5547000
$b2[11,1] = 5092000
$b2[11,2] = 5131000
$b2[11,3] = 5145000
$b2[11,4] = 4714000
$b2[11,5] = 4746000
$b2[11,6] = 4506000
$b2[11,7] = 4549000
$b2[11,8] = 4570000
$b2[12,0] = 0
$b2[12,1] = 0
$b2[12,2] = 0
$b2[12,3] = 0
$b2[12,4] = 0
$b2[12,5] = 0
$b2[12,6] = 0
$b2[12,7] = 0
$b2[12,8] = 0
$b2[13,0] = 0
$b2[13,1] = 0
$b2[13,2] = 0
$b2[13,3] = 0
$b2[13,4] = 0
$b2[13,5] = 0
$b2[13,6] = 0
$b2[13,7] = 0
$b2[13,8] = 0
$b2[14,0] = 1869000
$b2[14,1] = 1964000
$b2[14,2] = 1758000
$b2[14,3] = 1474000
$b2[14,4] = 923000
$b2[14,5] = 953000
$b2[14,6] = 964000
$b2[14,7] = 710000
$b2[14,8] = 694000
$b2[15,0] = 0
$b2[15,1] = 0
$b2[15,2] = 0
$b2[15,3] = 0
$b2[15,4] = 0
$b2[15,5] = 0
$b2[15,6] = 0
$b2[15,7] = 0
$b2[15,8] = 0
$b2[16,0] = 23962000
$b2[16,1] = 23409000
$b2[16,2] = 22088000
$b2[16,3] = 21206000
$b2[16,4] = 19594000
$b2[16,5] = 17175000
$b2[16,6] = 16241000
$b2[16,7] = 16643000
$b2[16,8] = 15869000
$b2[17,0] = $null
$b2[17,1] = $null
$b2[17,2] = $null
$b2[17,3] = $null
$b2[17,4] = $null
$b2[17,5] = $null
$b2[17,6] = $null
$b2[17,7] = $null
$b2[17,8] = $null
$b2[18,0] = $null
$b2[18,1] = $null
$b2[18,2] = $null
$b2[18,3] = $null
$b2[18,4] = $null
$b2[18,5] = $null
$b2[18,6] = $null
$b2[18,7] = $null
$b2[18,8] = $null
$b2[19,0] = 158000
$b2[19,1] = 119000
$b2[19,2] = 126000
$b2[19,3] = 15000
$b2[19,4] = 99000
$b2[19,5] = 116000
$b2[19,6] = 116000
$b2[19,7] = 125000
$b2[19,8] = 104000
$b2[20,0] = "NA"
$b2[20,1] = "NA"
$b2[20,2] = "NA"
$b2[20,3] = "NA"
$b2[20,4] = "NA"
$b2[20,5] = 680000
$b2[20,6] = 680000
$b2[20,7] = "NA"
$b2[20,8] = "NA"
$b2[21,0] = 4912000
$b2[21,1] = 4944000
$b2[21,2] = 4511000
$b2[21,3] = 4795000
$b2[21,4] = 4426000
$b2[21,5] = 4480000
$b2[21,6] = 4214000
$b2[21,7] = 4429000
$b2[21,8] = 3928000
$b2[22,0] = 5070000
$b2[22,1] = 5063000
$b2[22,2] = 4637000
$b2[22,3] = 4810000
$b2[22,4] = 4525000
$b2[22,5] = 5276000
$b2[22,6] = 5010000
$b2[22,7] = 4554000
$b2[22,8] = 4032000
$b2[23,0] = 4240000
$b2[23,1] = 4238000
$b2[23,2] = 4236000
$b2[23,3] = 4234000
$b2[23,4] = 4232000
$b2[23,5] = 820000
$b2[23,6] = 820000
$b2[23,7] = 1500000
$b2[23,8] = 1500000
$b2[24,0] = 3741000
$b2[24,1] = 3716000
$b2[24,2] = 3555000
$b2[24,3] = 3538000
$b2[24,4] = 2563000
$b2[24,5] = 2480000
$b2[24,6] = 2343000
$b2[24,7] = 2492000
$b2[24,8] = 2203000
$b2[25,0] = 0
$b2[25,1] = 0
$b2[25,2] = 0
$b2[25,3] = 0
$b2[25,4] = 0
$b2[25,5] = 0
$b2[25,6] = 0
$b2[25,7] = 0
$b2[25,8] = 0
$b2[26,0] = 0
$b2[26,1] = 0
$b2[26,2] = 0
$b2[26,3] = 0
$b2[26,4] = 0
$b2[26,5] = 0
$b2[26,6] = 0
$b2[26,7] = 0
$b2[26,8] = 0
$b2[27,0] = 0
$b2[27,1] = 0
$b2[27,2] = 0
$b2[27,3] = 0
$b2[27,4] = 0
$b2[27,5] = 0
$b2[27,6] = 0
$b2[27,7] = 0
$b2[27,8] = 0
$b2[28,0] = 13051000
$b2[28,1] = 13017000
$b2[28,2] = 12428000
$b2[28,3] = 12582000
$b2[28,4] = 11320000
$b2[28,5] = 8576000
$b2[28,6] = 8173000
$b2[28,7] = 8546000
$b2[28,8] = 7735000
$b2[29,0] = $null
$b2[29,1] = $null
$b2[29,2] = $null
$b2[29,3] = $null
$b2[29,4] = $null
$b2[29,5] = $null
$b2[29,6] = $null
$b2[29,7] = $null
$b2[29,8] = $null
$b2[30,0] = 0
$b2[30,1] = 0
$b2[30,2] = 0
$b2[30,3] = 0
$b2[30,4] = 0
$b2[30,5] = 0
$b2[30,6] = 0
$b2[30,7] = 0
$b2[30,8] = 0
$b2[31,0] = 0
$b2[31,1] = 0
$b2[31,2] = 0
$b2[31,3] = 0
$b2[31,4] = 0
$b2[31,5] = 0
$b2[31,6] = 0
$b2[31,7] = 0
$b2[31,8] = 0
$b2[32,0] = 0
$b2[32,1] = 0
$b2[32,2] = 0
$b2[32,3] = 0
$b2[32,4] = 0
$b2[32,5] = 0
$b2[32,6] = 0
$b2[32,7] = 0
$b2[32,8] = 0
$b2[33,0] = 0
$b2[33,1] = 0
$b2[33,2] = 0
$b2[33,3] = 0
$b2[33,4] = 0
$b2[33,5] = 0
$b2[33,6] = 0
$b2[33,7] = 0
$b2[33,8] = 0
$b2[34,0] = 9690000
$b2[34,1] = 9362000
$b2[34,2] = 8718000
$b2[34,3] = 7791000
$b2[34,4] = 7382000
$b2[34,5] = 6939000
$b2[34,6] = 6605000
$b2[34,7] = 6381000
$b2[34,8] = 5940000
$b2[35,0] = 0
$b2[35,1] = 0
$b2[35,2] = 0
$b2[35,3] = 0
$b2[35,4] = 0
$b2[35,5] = 0
$b2[35,6] = 0
$b2[35,7] = 0
$b2[35,8] = 0
$b2[36,0] = 0
$b2[36,1] = 0
$b2[36,2] = 0
$b2[36,3] = 0
$b2[36,4] = 0
$b2[36,5] = 0
$b2[36,6] = 0
$b2[36,7] = 0
$b2[36,8] = 0
$b2[37,0] = 0
$b2[37,1] = 0
$b2[37,2] = 0
$b2[37,3] = 0
$b2[37,4] = 0
$b2[37,5] = 0
$b2[37,6] = 0
$b2[37,7] = 0
$b2[37,8] = 0
$b2[38,0] = 10911000
$b2[38,1] = 10392000
$b2[38,2] = 9660000
$b2[38,3] = 8624000
$b2[38,4] = 8274000
$b2[38,5] = 8599000
$b2[38,6] = 8068000
$b2[38,7] = 8097000
$b2[38,8] = 8134000
$b2[39,0] = 0
$b2[39,1] = 0
$b2[39,2] = 0
$b2[39,3] = 0
$b2[39,4] = 0
$b2[39,5] = 0
$b2[39,6] = 0
$b2[39,7] = 0
$b2[39,8] = 0
$ws.Range("D38:L77").Value = $b2

$b3 = New-Object 'object[,]' 23,9
$b3[0,0] = 43406
$b3[0,1] = 43315
$b3[0,2] = 43224
$b3[0,3] = 43133
$b3[0,4] = 43042
$b3[0,5] = 42951
$b3[0,6] = 42860
$b3[0,7] = 42735
$b3[0,8] = 42643
$b3[1,0] = 334000
$b3[1,1] = 644000
$b3[1,2] = 942000
$b3[1,3] = -387000
$b3[1,4] = 395000
$b3[1,5] = 406000
$b3[1,6] = 245000
$b3[1,7] = 441000
$b3[1,8] = 319000
$b3[2,0] = $null
$b3[2,1] = $null
$b3[2,2] = $null
$b3[2,3] = $null
$b3[2,4] = $null
$b3[2,5] = $null
$b3[2,6] = $null
$b3[2,7] = $null
$b3[2,8] = $null
$b3[3,0] = 158000
$b3[3,1] = 151000
$b3[3,2] = 156000
$b3[3,3] = 146000
$b3[3,4] = 149000
$b3[3,5] = 134000
$b3[3,6] = 136000
$b3[3,7] = 84000
$b3[3,8] = 87000
$b3[4,0] = 0
$b3[4,1] = 0
$b3[4,2] = 0
$b3[4,3] = 0
$b3[4,4] = 0
$b3[4,5] = 0
$b3[4,6] = 0
$b3[4,7] = 0
$b3[4,8] = 0
$b3[5,0] = 0
$b3[5,1] = 0
$b3[5,2] = 0
$b3[5,3] = 0
$b3[5,4] = 0
$b3[5,5] = 0
$b3[5,6] = 0
$b3[5,7] = 0
$b3[5,8] = 0
$b3[6,0] = 0
$b3[6,1] = 0
$b3[6,2] = 0
$b3[6,3] = 0
$b3[6,4] = 0
$b3[6,5] = 0
$b3[6,6] = 0
$b3[6,7] = 0
$b3[6,8] = 0
$b3[7,0] = 0
$b3[7,1] = 0
$b3[7,2] = 0
$b3[7,3] = 0
$b3[7,4] = 0
$b3[7,5] = 0
$b3[7,6] = 0
$b3[7,7] = 0
$b3[7,8] = 0
$b3[8,0] = 0
$b3[8,1] = 0
$b3[8,2] = 0
$b3[8,3] = 0
$b3[8,4] = 0
$b3[8,5] = 0
$b3[8,6] = 0
$b3[8,7] = 0
$b3[8,8] = 0
$b3[9,0] = 769000
$b3[9,1] = 787000
$b3[9,2] = 1095000
$b3[9,3] = 851000
$b3[9,4] = 969000
$b3[9,5] = 621000
$b3[9,6] = 777000
$b3[9,7] = 464000
$b3[9,8] = 620000
$b3[10,0] = $null
$b3[10,1] = $null
$b3[10,2] = $null
$b3[10,3] = $null
$b3[10,4] = $null
$b3[10,5] = $null
$b3[10,6] = $null
$b3[10,7] = $null
$b3[10,8] = $null
$b3[11,0] = -57000
$b3[11,1] = -60000
$b3[11,2] = -61000
$b3[11,3] = -99000
$b3[11,4] = -59000
$b3[11,5] = -56000
$b3[11,6] = -49000
$b3[11,7] = -44000
$b3[11,8] = -30000
$b3[12,0] = 0
$b3[12,1] = 0
$b3[12,2] = 0
$b3[12,3] = 0
$b3[12,4] = 0
$b3[12,5] = 0
$b3[12,6] = 0
$b3[12,7] = 0
$b3[12,8] = 0
$b3[13,0] = 0
$b3[13,1] = 0
$b3[13,2] = 0
$b3[13,3] = 0
$b3[13,4] = 0
$b3[13,5] = 0
$b3[13,6] = 0
$b3[13,7] = 0
$b3[13,8] = 0
$b3[14,0] = 273000
$b3[14,1] = 322000
$b3[14,2] = 39000
$b3[14,3] = -651000
$b3[14,4] = -324000
$b3[14,5] = -944000
$b3[14,6] = 407000
$b3[14,7] = 287000
$b3[14,8] = 24000
$b3[15,0] = $null
$b3[15,1] = $null
$b3[15,2] = $null
$b3[15,3] = $null
$b3[15,4] = $null
$b3[15,5] = $null
$b3[15,6] = $null
$b3[15,7] = $null
$b3[15,8] = $null
$b3[16,0] = 0
$b3[16,1] = 0
$b3[16,2] = 0
$b3[16,3] = 0
$b3[16,4] = 0
$b3[16,5] = 0
$b3[16,6] = 0
$b3[16,7] = 0
$b3[16,8] = 0
$b3[17,0] = 0
$b3[17,1] = 0
$b3[17,2] = 0
$b3[17,3] = 0
$b3[17,4] = 0
$b3[17,5] = 0
$b3[17,6] = 0
$b3[17,7] = 0
$b3[17,8] = 0
$b3[18,0] = 0
$b3[18,1] = 0
$b3[18,2] = 0
$b3[18,3] = 0
$b3[18,4] = 0
$b3[18,5] = 0
$b3[18,6] = 0
$b3[18,7] = 0
$b3[18,8] = 0
$b3[19,0] = 0
$b3[19,1] = 0
$b3[19,2] = 0
$b3[19,3] = 0
$b3[19,4] = 0
$b3[19,5] = 0
$b3[19,6] = 0
$b3[19,7] = 0
$b3[19,8] = 0
$b3[20,0] = 37000
$b3[20,1] = -89000
$b3[20,2] = -3000
$b3[20,3] = -231000
$b3[20,4] = 1815000
$b3[20,5] = 12000
$b3[20,6] = -538000
$b3[20,7] = -615000
$b3[20,8] = -981000
$b3[21,0] = 0
$b3[21,1] = 0
$b3[21,2] = 0
$b3[21,3] = 0
$b3[21,4] = 0
$b3[21,5] = 0
$b3[21,6] = 0
$b3[21,7] = 0
$b3[21,8] = 0
$b3[22,0] = 1079000
$b3[22,1] = 1020000
$b3[22,2] = 1131000
$b3[22,3] = -31000
$b3[22,4] = 2460000
$b3[22,5] = -309000
$b3[22,6] = 644000
$b3[22,7] = 136000
$b3[22,8] = -337000
$ws.Range("D80:L102").Value = $b3
